$wb = $excel.ActiveWorkbook

# --- Sheet: Collection ---
$wsCollection = $wb.Worksheets.Item("Collection")
$wsCollection.Range("B3").Value = "https://dx.doi.org/10.5281/zenodo.1009240"
$wsCollection.Hyperlinks.Add($wsCollection.Range("B3"), "https://dx.doi.org/10.5281/zenodo.1009240")

# --- Sheet: People ---
$wsPeople = $wb.Worksheets.Item("People")
$wsPeople.Range("C1").Value = "givenName"
$wsPeople.Range("D1").Value = "familyName"

# --- Sheet: Licenses ---
$wsLicenses = $wb.Worksheets.Item("Licenses")
$wsLicenses.Range("C2").Value = "This work is licensed under the Creative Commons Attribution-NonCommercial-ShareAlike 3.0 Australia License. To view a copy of this license, visit http://creativecommons.org/licenses/by-nc-sa/3.0/au/ or send a letter to Creative Commons, PO Box 1866, Mountain View, CA 94042, USA."
$wsLicenses.Range("D2").Value = "CreativeWork"
